# Applies the cryptos.xlsx price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.892.32"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.039.57"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'227.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'60.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.41%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("D10").Value = "'0.0812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").Value = "'14.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "2.349.37"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").Value = "'20.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").Value = "'0.754"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'5.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "2.060.20"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "37.889.87"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "'69.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "0.0₃0824"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "'223.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "'166.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "'9.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("D28").Value = "'0.130"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").Value = "'18.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.40%  "
$ws.Range("D30").Value = "'1.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").Value = "'4.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").Value = "'4.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.31%  "
$ws.Range("D35").Value = "'0.0602"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "'6.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.72%  "
$ws.Range("D37").Value = "'2.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.07%  "
$ws.Range("D38").Value = "'3.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "1.535.08"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'96.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").Value = "'2.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'16.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.53%  "
$ws.Range("D45").Value = "'0.0923"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "'1.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'3.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.25%  "
$ws.Range("D48").Value = "'2.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'7.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").Value = "2.239.74"
$ws.Range("E51").Value = "  +0.54%  "
